$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and week-covering dates) ---
$ws.Range("A8").Value = "Volume 29   Number  44"
$ws.Range("C9").Value = "Report Covering the Week  10/31/2022  Through  11/6/2022"

# --- Text-literal "0" / "***.*" cells: copy format+value from an existing donor cell ---
# C15 already holds text "0" (style s=14); E15 already holds text "***.*" (style s=14)
$ws.Range("C15").Copy($ws.Range("F15"))
$ws.Range("C15").Copy($ws.Range("D16"))
$ws.Range("E15").Copy($ws.Range("E16"))
$ws.Range("C15").Copy($ws.Range("C18"))
$ws.Range("C15").Copy($ws.Range("D18"))
$ws.Range("E15").Copy($ws.Range("E18"))
$ws.Range("C15").Copy($ws.Range("D22"))
$ws.Range("E15").Copy($ws.Range("E22"))
$ws.Range("C15").Copy($ws.Range("C25"))
$ws.Range("C15").Copy($ws.Range("F26"))
$ws.Range("C15").Copy($ws.Range("D27"))
$ws.Range("E15").Copy($ws.Range("E27"))

# --- C27 flips from text "0" back to a plain number; pull the numeric
# column's format (from C16, already style 15) before writing the value ---
$ws.Range("C16").Copy()
$ws.Range("C27").PasteSpecial(-4122)

# --- Plain numeric updates ---
$ws.Range("N15").Value = -34.782608695652
$ws.Range("C16").Value = 4
$ws.Range("F16").Value = 13
$ws.Range("G16").Value = 7
$ws.Range("H16").Value = 85.714285714285
$ws.Range("I16").Value = 94
$ws.Range("K16").Value = 54.098360655737
$ws.Range("L16").Value = 11.904761904761
$ws.Range("M16").Value = -33.802816901408
$ws.Range("N16").Value = -85.970149253731
$ws.Range("D17").Value = 1
$ws.Range("E17").Value = 0
$ws.Range("F17").Value = 4
$ws.Range("G17").Value = 24
$ws.Range("H17").Value = -83.333333333333
$ws.Range("I17").Value = 124
$ws.Range("J17").Value = 137
$ws.Range("K17").Value = -9.489051094890
$ws.Range("L17").Value = 20.388349514563
$ws.Range("M17").Value = 30.526315789473
$ws.Range("N17").Value = -51.181102362204
$ws.Range("F18").Value = 5
$ws.Range("G18").Value = 11
$ws.Range("H18").Value = -54.545454545454
$ws.Range("L18").Value = 18.604651162790
$ws.Range("M18").Value = 56.923076923076
$ws.Range("N18").Value = -80
$ws.Range("C19").Value = 8
$ws.Range("D19").Value = 6
$ws.Range("E19").Value = 33.333333333333
$ws.Range("F19").Value = 26
$ws.Range("G19").Value = 32
$ws.Range("H19").Value = -18.75
$ws.Range("I19").Value = 292
$ws.Range("J19").Value = 224
$ws.Range("K19").Value = 30.357142857142
$ws.Range("L19").Value = 36.448598130841
$ws.Range("M19").Value = 23.206751054852
$ws.Range("N19").Value = -49.305555555555
$ws.Range("C20").Value = 1
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 9
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = 80
$ws.Range("I20").Value = 67
$ws.Range("J20").Value = 44
$ws.Range("K20").Value = 52.272727272727
$ws.Range("L20").Value = 67.5
$ws.Range("M20").Value = 179.166666666667
$ws.Range("N20").Value = -79.0625
$ws.Range("D21").Value = 8
$ws.Range("E21").Value = 75
$ws.Range("F21").Value = 57
$ws.Range("G21").Value = 79
$ws.Range("H21").Value = -27.848101265822
$ws.Range("I21").Value = 695
$ws.Range("J21").Value = 526
$ws.Range("K21").Value = 32.129277566539
$ws.Range("L21").Value = 28.942486085343
$ws.Range("M21").Value = 20.450606585788
$ws.Range("N21").Value = -70.600676818950
$ws.Range("G22").Value = 2
$ws.Range("H22").Value = 50
$ws.Range("I22").Value = 18
$ws.Range("K22").Value = -28
$ws.Range("L22").Value = -28
$ws.Range("M22").Value = -35.714285714285
$ws.Range("C23").Value = 2
$ws.Range("D23").Value = 1
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 7
$ws.Range("H23").Value = -12.5
$ws.Range("I23").Value = 98
$ws.Range("J23").Value = 93
$ws.Range("K23").Value = 5.376344086021
$ws.Range("L23").Value = -3.921568627450
$ws.Range("M23").Value = 44.117647058823
$ws.Range("C24").Value = 14
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 55.555555555555
$ws.Range("F24").Value = 38
$ws.Range("G24").Value = 52
$ws.Range("H24").Value = -26.923076923076
$ws.Range("I24").Value = 443
$ws.Range("J24").Value = 365
$ws.Range("K24").Value = 21.369863013698
$ws.Range("L24").Value = -21.033868092691
$ws.Range("M24").Value = -17.041198501872
$ws.Range("D25").Value = 3
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 12
$ws.Range("G25").Value = 19
$ws.Range("H25").Value = -36.842105263157
$ws.Range("J25").Value = 177
$ws.Range("K25").Value = 2.259887005649
$ws.Range("L25").Value = -3.723404255319
$ws.Range("M25").Value = -29.296875
$ws.Range("H26").Value = -100
$ws.Range("C27").Value = 1
$ws.Range("F27").Value = 8
$ws.Range("G27").Value = 3
$ws.Range("H27").Value = 166.666666666667
$ws.Range("I27").Value = 47
$ws.Range("K27").Value = 34.285714285714
$ws.Range("L27").Value = 88
